# Update the "想去人数" (wanted-to-go count) figures on the "展览" and
# "全部类型" sheets to reflect the newly generated data.
#   F2: 526 -> 527
#   F4: 12  -> 14

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 527
    $ws.Range("F4").Value = 14
}
